$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "98.471.61"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "3.432.47"
$ws.Range("E3").Value = "  +2.78%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'258.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.17%  "
$ws.Range("D6").Value = "'660.52"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.85%  "
$ws.Range("D7").Value = "'1.49"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.01%  "
$ws.Range("D8").Value = "'0.446"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.10%  "
$ws.Range("E9").Value = "  -0.41%  "
$ws.Range("E10").Value = "  +0.01%  "
$ws.Range("D11").Value = "3.430.96"
$ws.Range("E11").Value = "  +2.84%  "
$ws.Range("E12").Value = "  +4.25%  "
$ws.Range("D13").Value = "'42.85"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.41%  "
$ws.Range("D14").Value = "'6.40"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +15.92%  "
$ws.Range("D15").Value = "'0.0000272"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.74%  "
$ws.Range("D16").Value = "98.139.85"
$ws.Range("E16").Value = "  -0.55%  "
$ws.Range("D17").Value = "4.071.96"
$ws.Range("E17").Value = "  +3.15%  "
$ws.Range("D18").Value = "'9.35"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +32.39%  "
$ws.Range("B19").Value = "Stellar"
$ws.Range("C19").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D19").Value = "'0.613"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +39.00%  "
$ws.Range("B20").Value = "WrappedEther"
$ws.Range("C20").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D20").Value = "3.431.10"
$ws.Range("E20").Value = "  +2.96%  "
$ws.Range("D21").Value = "'17.96"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.66%  "
$ws.Range("D22").Value = "'3.53"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.84%  "
$ws.Range("D23").Value = "'10.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.56%  "
$ws.Range("D24").Value = "'521.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.53%  "
$ws.Range("D25").Value = "'0.0000209"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.38%  "
$ws.Range("D26").Value = "'6.41"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.67%  "
$ws.Range("D27").Value = "'102.68"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.60%  "
$ws.Range("D28").Value = "'13.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.04%  "
$ws.Range("D29").Value = "3.616.40"
$ws.Range("E29").Value = "  +3.03%  "
$ws.Range("D30").Value = "'0.157"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.58%  "
$ws.Range("E31").Value = "  +8.74%  "
$ws.Range("D32").Value = "'0.201"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.83%  "
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("D34").Value = "'0.593"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +13.21%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("B36").Value = "PancakeSwap"
$ws.Range("C36").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D36").Value = "'2.36"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +15.47%  "
$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D37").Value = "'30.31"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.40%  "
$ws.Range("E38").Value = "  +4.87%  "
$ws.Range("D39").Value = "'1.46"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +12.58%  "
$ws.Range("D40").Value = "'535.31"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.81%  "
$ws.Range("E41").Value = "  +0.68%  "
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("B43").Value = "Cosmos"
$ws.Range("C43").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D43").Value = "'9.29"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +19.43%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "'0.0447"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +13.05%  "
$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").Value = "'0.886"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +9.15%  "
$ws.Range("B46").Value = "WhiteBITCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D46").Value = "'24.78"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.11%  "
$ws.Range("D47").Value = "'5.96"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +21.32%  "
$ws.Range("D48").Value = "'3.73"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.07%  "
$ws.Range("D49").Value = "'3.36"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.83%  "
$ws.Range("D50").Value = "'1.69"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +14.96%  "
$ws.Range("E51").Value = "  +4.58%  "
